$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "Right-click the layer in the layer pane" paragraph: append a
#    clarification about groups.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("ck the layer in the layer pane", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)   # wdCollapseEnd
    $rng.InsertAfter(" (or on the group, if there are several image layers in a group)")
}

# ------------------------------------------------------------------
# 2. "Duplicate the layer " paragraph: drop the trailing space from
#    the plain-text run ("uplicate the layer " -> "uplicate the layer")
#    and add " (or the group)" right before the bold "as a new file".
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("uplicate the layer ", $true, $false, $false, $false, $false, $true, 1, $false, "uplicate the layer", 2) | Out-Null

$rng = $d.Content
$found = $rng.Find.Execute("as a new file", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.InsertBefore(" (or the group)")
}

# ------------------------------------------------------------------
# 3. Relocate the "_GoBack" bookmark from the end of the document to
#    right before the bold "as a new file" run, followed by a single
#    space run.
# ------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$rng = $d.Content
$found = $rng.Find.Execute("as a new file", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $bmRange = $d.Range($rng.Start, $rng.Start)
    $d.Bookmarks.Add("_GoBack", $bmRange)
    $rng.InsertBefore(" ")
}
